$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows 1-16 shift down to 2-17.
$ws.Rows("1:1").Insert()

# Header labels (order_id, status, create_time, express_company); E1 stays blank.
$ws.Range("A1").Value = "order_id"
$ws.Range("B1").Value = "status"
$ws.Range("C1").Value = "create_time"
$ws.Range("D1").Value = "express_company"

# Header row formatting: 12pt Verdana (distinct font entry from the body's 9pt Verdana).
$ws.Range("A1:E1").Font.Size = 12
$ws.Range("A1:E1").Font.Name = "Verdana"

# Fill in the previously-missing status/create_time/express_company/extra columns
# for the data rows (now at rows 2-17), following the existing sequence.
$data = @(
    @(2, 22, 43429, 233, 322),
    @(3, 23, 43430, 234, 323),
    @(4, 24, 43431, 235, 323),
    @(5, 25, 43432, 236, 323),
    @(6, 26, 43433, 237, 323),
    @(7, 27, 43434, 238, 323),
    @(8, 28, 43435, 239, 323),
    @(9, 29, 43436, 240, 323),
    @(10, 30, 43437, 241, 323),
    @(11, 31, 43438, 242, 323),
    @(12, 32, 43439, 243, 323),
    @(13, 33, 43440, 244, 323),
    @(14, 34, 43441, 245, 323),
    @(15, 35, 43442, 246, 323),
    @(16, 36, 43443, 247, 323),
    @(17, 37, 43444, 248, 323)
)

foreach ($row in $data) {
    $r = $row[0]
    if ($r -gt 2) {
        # row 2 (the original row 1) already carries the date format on C2;
        # the rest need the C column's date style copied over before writing.
        $ws.Range("C2").Copy()
        $ws.Range("C" + $r).PasteSpecial(-4122)
    }
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]
    $ws.Range("E" + $r).Value = $row[4]
}

# New row 17's order_id follows the same pattern as the other data rows.
$ws.Range("A17").Value = 944138813511303000

# Restore the selection to match the edited workbook.
$ws.Range("F7").Select()
